# EntendiendoBN.xlsx update: add gamma/beta scale+shift (Norma 1) and
# per-sample variance columns to the batch-norm worksheet, shifting the
# existing Voltaje X / Normalizados table one column to the right and
# zeroing out Epsilon.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Keep a border-styled source cell (B1 already has the thin-box style used
# throughout the sheet) to stamp onto the newly occupied cells before we
# wipe the old layout.
$ws.Range("B1").Copy()
$ws.Range("C2:D13").PasteSpecial(-4122)
$ws.Range("H2:H14").PasteSpecial(-4122)
$ws.Range("B15:C17").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---- wipe the old (pre-shift) layout that's no longer part of the table ----
$ws.Range("B1:D14").ClearContents()
$ws.Range("A14:A16").ClearContents()
$ws.Range("B19:D19").ClearContents()

# ---- headers (row 2) ----
$ws.Range("C2").Value = "Voltaje X"
$ws.Range("D2").Value = "Normalizados"
$ws.Range("F2").Value = "Norma 1"
$ws.Range("H2").Value = "VAR"

# ---- data rows 3-13 (Voltaje X + Normalizados + Norma 1 + VAR) ----
$xvals = @(1, 1.5, 2, 2.5, 3, 3.5, 4, 4.5, 5, 5.5, 6)
for ($i = 0; $i -lt $xvals.Length; $i++) {
    $r = 3 + $i
    $ws.Cells.Item($r, 3).Value = $xvals[$i]
    $ws.Cells.Item($r, 4).Formula = "=(C$r-`$C`$15)/(SQRT(`$C`$16+`$C`$17))"
    $ws.Cells.Item($r, 6).Formula = "=`$C`$19*(C$r-`$C`$15)/(SQRT(`$C`$16+`$C`$17))+`$C`$20"
    $ws.Cells.Item($r, 8).Formula = "=(C$r-`$C`$15)^2"
}

# ---- VAR average ----
$ws.Range("H14").Formula = "=SUM(H3:H13)/11"

# ---- Media / Varianza / Epsilon ----
$ws.Range("B15").Value = "Media"
$ws.Range("C15").Formula = "=AVERAGE(C3:C13)"
$ws.Range("B16").Value = "Varianza"
$ws.Range("C16").Value = 2.5
$ws.Range("B17").Value = "Epsilon"
$ws.Range("C17").Value = 0

# ---- Gamma / Beta ----
$ws.Range("B19").Value = "Gamma"
$ws.Range("C19").Value = 0.6582914
$ws.Range("B20").Value = "Beta"
$ws.Range("C20").Value = 0.27094355

# numeric formatting for the Norma 1 column
$ws.Range("F3:F13").NumberFormat = "0.0000000"

# column widths
$ws.Columns.Item(3).ColumnWidth = 11.85546875
$ws.Columns.Item(4).ColumnWidth = 13.140625
$ws.Columns.Item(6).ColumnWidth = 12.28515625

# selection
$ws.Range("D13:E13").Select()
